$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the closing-border formatting (row 23, last row of the data table)
#    onto row 21, which will become the new last row once rows 22-23 are
#    removed below.
$ws.Range("B23:J23").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) Remove the two obsolete worker rows (old rows 22 and 23); this shifts
#    every row below them up by two (signature block 28/29 -> 26/27).
$ws.Rows("22:23").Delete()

# 3) Write the new worker roster into the (now 6-row) data table.
$ws.Cells.Item(16, 3).Value = "1128056761"
$ws.Cells.Item(16, 4).Value = "DEIVIS DAVID ASENCIO ORTIZ"
$ws.Cells.Item(16, 5).Value = "2507"
$ws.Cells.Item(16, 6).Value = 64000
$ws.Cells.Item(16, 7).Value = 1600000

$ws.Cells.Item(17, 3).Value = "1143371014"
$ws.Cells.Item(17, 4).Value = "MARIA FERNANDA BALCEIRO MORANTE"
$ws.Cells.Item(17, 5).Value = "2402"
$ws.Cells.Item(17, 6).Value = 52000
$ws.Cells.Item(17, 7).Value = 1200000

$ws.Cells.Item(18, 3).Value = "1047386048"
$ws.Cells.Item(18, 4).Value = "CRISTIAN DAVID OCHOA ALVEAR"
$ws.Cells.Item(18, 5).Value = "2402"
$ws.Cells.Item(18, 6).Value = 52000
$ws.Cells.Item(18, 7).Value = 1300000

$ws.Cells.Item(19, 3).Value = "1002195062"
$ws.Cells.Item(19, 4).Value = "DIEGO ANDRES HERNANDEZ CABARCAS"
$ws.Cells.Item(19, 5).Value = "2507"
$ws.Cells.Item(19, 6).Value = 56940
$ws.Cells.Item(19, 7).Value = 1423500

$ws.Cells.Item(20, 3).Value = "1002244060"
$ws.Cells.Item(20, 4).Value = "CHRISTIAN ARIEL BALCEIRO MORANTE"
$ws.Cells.Item(20, 5).Value = "2507"
$ws.Cells.Item(20, 6).Value = 74000
$ws.Cells.Item(20, 7).Value = 1850000

$ws.Cells.Item(21, 3).Value = "1049454782"
$ws.Cells.Item(21, 4).Value = "JAIME ANDRES QUIROZ GUERRERO"
$ws.Cells.Item(21, 5).Value = "2507"
$ws.Cells.Item(21, 6).Value = 56940
$ws.Cells.Item(21, 7).Value = 1423500

# 4) Refresh the "Valor Mora" total to match the new roster.
$ws.Range("E11").Value = 355880

